# Registration.xlsx update: add "DuplicateRegis" and "MandatoryFields" sheets,
# add new Thank-you/Confirmation/Continue/Welcome columns to Register, swap
# the sample applicant name + contact email, and add the duplicate-account
# error message used on the two new sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Register"

# ------------------------------------------------------------------
# Register (sheet1): extend the header/sample row with four new fields
# ------------------------------------------------------------------

# Preserve the old I1/I2 ("runMode"/"Y") over in the new last column (M)
$ws1.Range("M1").Value = "runMode"
$ws1.Range("M2").Value = "Y"

# New header cells (row 1)
$ws1.Range("I1").Value = "ThankyouText"

$ws1.Range("J1").NumberFormat = "@"
$ws1.Range("J1").Value = "ConfirmationText"

$ws1.Range("K1").NumberFormat = "@"
$ws1.Range("K1").Value = "ContinueButtonText"

$ws1.Range("L1").NumberFormat = "@"
$ws1.Range("L1").Value = "WelcomeText"

# New sample-value cells (row 2)
$ws1.Range("I2").NumberFormat = "@"
$ws1.Range("I2").Value = "Thank you for creating a user account."

$ws1.Range("J2").NumberFormat = "@"
$ws1.Range("J2").Value = "Confirmation has been sent to your email address."

$ws1.Range("K2").NumberFormat = "@"
$ws1.Range("K2").Value = "Continue with application"

$ws1.Range("L2").NumberFormat = "@"
$ws1.Range("L2").Value = "Welcome Howard"

# Swap the sample applicant's name
$ws1.Range("A2").Value = "Grant"
$ws1.Range("B2").Value = "Howard"

# Update the contact e-mail (value + underlying mailto hyperlink)
$ws1.Range("E2").Hyperlinks.Delete()
$ws1.Range("E2").Value = "pankaj.missguided1354@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("E2"), "mailto:pankaj.missguided1354@gmail.com")
$ws1.Range("E2").NumberFormat = "@"

# Approximate column widths for the new columns (cosmetic only)
$ws1.Columns.Item(9).ColumnWidth = 32.66
$ws1.Columns.Item(10).ColumnWidth = 22.89
$ws1.Columns.Item(11).ColumnWidth = 24.11
$ws1.Columns.Item(12).ColumnWidth = 24.11

$ws1.Rows("1:2").Select()

# ------------------------------------------------------------------
# DuplicateRegis (sheet2): copy of Register plus an ErrorMsg column
# ------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "DuplicateRegis"

$ws2.Range("A1").Value = "FamilyName"
$ws2.Range("B1").Value = "GivenName"
$ws2.Range("C1").Value = "DateOFBirth"
$ws2.Range("D1").Value = "Gender"
$ws2.Range("E1").Value = "EmailId"
$ws2.Range("F1").Value = "Nationality"
$ws2.Range("G1").Value = "Password"
$ws2.Range("H1").Value = "ConPassword"
$ws2.Range("I1").NumberFormat = "@"
$ws2.Range("I1").Value = "ErrorMsg"
$ws2.Range("J1").Value = "runMode"

$ws2.Range("A2").Value = "Grant"
$ws2.Range("B2").Value = "Howard"
$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Value = "'09-Oct-2005"
$ws2.Range("D2").Value = "M"
$ws2.Range("E2").Value = "pankaj.missguided1354@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("E2"), "mailto:pankaj.missguided1354@gmail.com")
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("F2").Value = "UK National"
$ws2.Range("G2").NumberFormat = "@"
$ws2.Range("G2").Value = "Pa55w0rd#"
$ws2.Range("H2").NumberFormat = "@"
$ws2.Range("H2").Value = "Pa55w0rd#"
$ws2.Range("I2").NumberFormat = "@"
$ws2.Range("I2").Value = "'A person with the same details already has a user account. If you know your username and password, select 'Log on' below. If you need assistance, select 'Forgot password?'"
$ws2.Range("J2").Value = "Y"

$ws2.Rows("1:2").Select()

# ------------------------------------------------------------------
# MandatoryFields (sheet3): minimal field list plus ErrorMsg
# ------------------------------------------------------------------

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "MandatoryFields"

$ws3.Range("A1").Value = "GivenName"
$ws3.Range("B1").Value = "Gender"
$ws3.Range("C1").NumberFormat = "@"
$ws3.Range("C1").Value = "ErrorMsg"
$ws3.Range("D1").Value = "runMode"

$ws3.Range("A2").Value = "Howard"
$ws3.Range("B2").Value = "M"
$ws3.Range("C2").NumberFormat = "@"
$ws3.Range("C2").Value = "'A person with the same details already has a user account. If you know your username and password, select 'Log on' below. If you need assistance, select 'Forgot password?'"
$ws3.Range("D2").Value = "Y"

$ws3.Columns.Item(3).Select()
$ws3.Activate()

Write-Output "Registration.xlsx updated: sheets, headers and sample rows set."
